# Update cryptocurrency price/volume data in-place (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "49.718.13"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +2.89%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.627.71"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +4.83%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "330.57"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.85%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "110.07"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.58%  "

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.535"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.560"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.39%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "40.87"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.78"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("E12").Value = "  +0.28%  "

# Row 13
$ws.Range("E13").Value = "  +0.78%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.30"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.038.74"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +4.84%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.615.20"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +4.31%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.877"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.81%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "49.702.00"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.21%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.11"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +11.32%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.39"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.83"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0955"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "281.82"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "72.81"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25
$ws.Range("E25").Value = "  +1.32%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "26.66"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.35%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("E28").Value = "  -2.69%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.99"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.92%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.144"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.10%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "36.27"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +2.67%  "

# Row 32
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "19.78"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "

# Row 34
$ws.Range("E34").Value = "  +1.89%  "

# Row 36
$ws.Range("E36").Value = "  +1.47%  "

# Row 38
$ws.Range("E38").Value = "  +2.17%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +5.81%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "22.84"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +5.15%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "123.32"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.53%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("E44").Value = "  +3.76%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "3.37"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +6.28%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.054.03"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.89%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +12.66%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.01"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +8.78%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.05"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

# Row 50
$ws.Range("E50").Value = "  +3.62%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "81.93"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.41%  "

